# "programação das figuras do grupo 16"
# Update the "Ano" (Year) column (C) data: the dates were stored as
# 31/12/<year> and need to become 01/01/<year> for every data row
# (rows 2-91), keeping the year itself unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$anoRange = $ws.Range("C2:C91")

# Force the cells to stay text (rather than letting Excel reinterpret
# "01/01/2016" style strings as real dates) while we rewrite the values.
$anoRange.NumberFormat = "@"
$anoRange.Replace("31/12/", "01/01/")

# Restore the default cell style so the cells end up identical in
# formatting to how they started (no explicit number format applied).
$anoRange.Style = "Normal"
